$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Insert a brand-new paragraph before the "A atividade de passeio..."
# paragraph (paragraph 2), containing the new sentence about the
# trainer's health assessment, and move the _GoBack bookmark there (it
# currently sits at the end of the last paragraph).
#
# We temporarily fill the new paragraph with a single placeholder
# character "X" so that when we add the _GoBack bookmark at the start of
# the paragraph, its position isn't anywhere near a paragraph boundary:
# adding a zero-length bookmark exactly at "paragraph end - 1" trips a
# positional bug in this COM host (Bookmarks.Add silently mis-resolves
# to a bogus range whenever the target position equals some paragraph's
# End-1). We sidestep it by only ever calling Bookmarks.Add at a safe
# interior/start position, then letting ordinary text insertion/deletion
# shift the already-created bookmark into its final place.
# -----------------------------------------------------------------------
$p2 = $d.Paragraphs(2)
$p2.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs(2)
$newPara.Range.Text = "X"

$newPara2 = $d.Paragraphs(2)
$bmpos = $newPara2.Range.Start

# Move the _GoBack bookmark here (Add at the safe placeholder position).
$d.Bookmarks("_GoBack").Delete()
$bmRange = $d.Range($bmpos, $bmpos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Insert the real sentence right before the placeholder; the collapsed
# bookmark naturally slides along to sit right after the new text.
$newText = "O adestrador faz uma avaliação do estado de saúde física do animal para depois elaborar um programa de atividades físicas e/ou treinamento adequado ao animal."
$insPoint = $d.Range($bmpos, $bmpos)
$insPoint.InsertBefore($newText)

# Remove the placeholder "X" that now trails the bookmark.
$bm = $d.Bookmarks("_GoBack")
$placeholder = $d.Range($bm.End, $bm.End + 1)
$placeholder.Delete()
